$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 417, shifting rows 417:515 down to 418:516.
$ws.Rows(417).Insert()

# Populate the newly inserted row 417 with the new record's data.
$ws.Range("A417").Value = 10
$ws.Range("B417").Value = "Vega Modelo de Temuco"
$ws.Range("C417").Value = "La Araucanía"
$ws.Range("D417").Value = 44932
$ws.Range("D417").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E417").Value = 9
$ws.Range("F417").Value = 100114014
$ws.Range("G417").Value = "Betarraga"
$ws.Range("H417").Value = "Sin especificar"
$ws.Range("I417").Value = "Primera"
$ws.Range("J417").Value = 40
$ws.Range("K417").Value = 9000
$ws.Range("L417").Value = 9000
$ws.Range("M417").Value = 9000
$ws.Range("N417").Value = '$/docena de paquetes'
$ws.Range("O417").Value = "Provincia de Cautín"
$ws.Range("P417").Value = 750
$ws.Range("Q417").Value = 12
$ws.Range("R417").Value = "Hortaliza"
